# Append new "S4 / Anderson" session rows (2013-06-03) to Sheet1 of the
# watchFftDataset workbook, matching the rest of the SSVEP/hybrid session log.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rows = @(
    @("S4","Anderson",41339,"2013-06-03-anderson","2013-06-03-14-44-02","ssvep-12Hz",12,0),
    @("S4","Anderson",41339,"2013-06-03-anderson","2013-06-03-14-51-37","hybrid-12Hz",12,1),
    @("S4","Anderson",41339,"2013-06-03-anderson","2013-06-03-14-59-20","ssvep-10Hz",10,0),
    @("S4","Anderson",41339,"2013-06-03-anderson","2013-06-03-15-05-14","ssvep-15Hz",15,0),
    @("S4","Anderson",41339,"2013-06-03-anderson","2013-06-03-15-11-19","hybrid-15Hz",15,1),
    @("S4","Anderson",41339,"2013-06-03-anderson","2013-06-03-15-17-38","hybrid-10Hz",10,1),
    @("S4","Anderson",41339,"2013-06-03-anderson","2013-06-03-15-23-51","ssvep-8-57Hz",8.57,0),
    @("S4","Anderson",41339,"2013-06-03-anderson","2013-06-03-15-29-39","hybrid-8-57Hz",8.57,1)
)

$startRow = 26
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]

    $ws.Cells.Item($r, 1).Value2 = $row[0]
    $ws.Cells.Item($r, 2).Value2 = $row[1]
    $ws.Cells.Item($r, 3).Value2 = $row[2]
    $ws.Cells.Item($r, 4).Value2 = $row[3]
    $ws.Cells.Item($r, 5).Value2 = $row[4]
    $ws.Cells.Item($r, 6).Value2 = $row[5]
    $ws.Cells.Item($r, 7).Value2 = $row[6]
    $ws.Cells.Item($r, 8).Value2 = $row[7]
}

# Copy the direct cell formatting (incl. the date number format in column C)
# from the previous data row down onto the newly added rows.
$ws.Range("A25:H25").Copy() | Out-Null
$ws.Range("A26:H33").PasteSpecial(-4122) | Out-Null

# Match the bestFit-style column width recalculation Excel performs once the
# new (wider) sessionDirectory / condition strings are present.
$ws.Columns.Item(4).ColumnWidth = 18.833333333333336
$ws.Columns.Item(6).ColumnWidth = 12.666666666666666

$ws.Range("R14").Select() | Out-Null
